$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.136.65'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -1.95%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.667.06'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -1.27%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('E5').Value = '  +0.00%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5118'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +3.47%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2635'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  +1.97%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06401'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  +5.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.61'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  -0.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07420'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  +1.82%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.675.82'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  -0.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.513'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +1.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5813'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  +1.78%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008576'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  +5.24%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.26'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -0.43%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.199.26'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  -1.79%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.927'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  -1.42%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.006'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  -0.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.83'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').Value = '  +1.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.89'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  +4.50%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.203'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  +0.74%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.007'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.40'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +0.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.631'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  +1.49%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1195'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  +6.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.60'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  +2.73%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06287'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  +12.89%  '

$ws.Range('E29').Value = '  -1.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.318'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  -0.19%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.528'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  +2.03%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.509'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  +1.94%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.640'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -0.05%  '

$ws.Range('E34').Value = '  +0.87%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6076'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  +4.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.364'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -1.59%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.648'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').Value = '  +1.34%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.166'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  +5.32%  '

$ws.Range('E39').Value = '  +1.71%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.080.48'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  +1.53%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8654'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  +2.19%  '

$ws.Range('E42').Value = '  +0.60%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.09'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  +3.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.816.23'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -1.69%  '

$ws.Range('E45').Value = '  +4.62%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.12'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  +0.02%  '

$ws.Range('E47').Value = '  +0.44%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.067'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  +0.28%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05201'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +0.25%  '

$ws.Range('E50').Value = '  -1.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.901'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  +6.10%  '
